$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Ativacao:" date text value update: 01/01/2019 -> 01/01/2023 ---
# The cells store the date as literal text (shared string), not a real date
# serial. A plain .Value assignment of a dd/mm/yyyy-looking string gets
# auto-parsed into a date serial by Excel, which would also disturb the
# cell's number format / style. To avoid that, stage the text as a text
# formula result in a scratch cell far away, copy it, and paste-special
# just the *values* into the target cells (keeps original formatting/style
# untouched, keeps the shared-string text type).
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""01/01/2023"""
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()

# --- 2) New English "Objectives" paragraph, mirrored into B11/C11 ---
# Row 11 currently only has A11 ("Objectives:"). Copy formatting from the
# analogous populated row (10) so the new cells pick up the same styles
# (B column style / C column "red" style) used throughout the sheet.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Value = "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."
$ws.Range("C11").Value = "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."

# --- 3) New English "Short syllabus" paragraph, mirrored into B13/C13 ---
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B13").Value = "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"
$ws.Range("C13").Value = "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"

# --- 4) New English "Syllabus" paragraph, mirrored into B15/C15 ---
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats
$syllabus = "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace" + [char]0x2019 + "s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."
$ws.Range("B15").Value = $syllabus
$ws.Range("C15").Value = $syllabus

$excel.CutCopyMode = $false
